# Horarios actualizados Linea 141 - 648
# Update scrape metadata + reorder/extend schedule rows across all 3 sheets
$wb = $excel.ActiveWorkbook

$wsLP1912 = $wb.Worksheets.Item("LP1912")
$wsLP215  = $wb.Worksheets.Item("LP1912-215")
$wsL6203  = $wb.Worksheets.Item("6203-6173")

# --- Sheet LP1912: header metadata ---
$wsLP1912.Range("A2").Value2 = "Última actualización: 14:24:16"
$wsLP1912.Range("A3").Value2 = "Total filas: 284"

# --- Sheet LP1912: swap mis-ordered duplicate-arrival rows ---
# swap rows 32 <-> 33
$tmpA = $wsLP1912.Cells.Item(32,1).Value2
$tmpC = $wsLP1912.Cells.Item(32,3).Value2
$tmpD = $wsLP1912.Cells.Item(32,4).Value2
$wsLP1912.Cells.Item(32,1).Value2 = $wsLP1912.Cells.Item(33,1).Value2
$wsLP1912.Cells.Item(32,3).Value2 = $wsLP1912.Cells.Item(33,3).Value2
$wsLP1912.Cells.Item(32,4).Value2 = $wsLP1912.Cells.Item(33,4).Value2
$wsLP1912.Cells.Item(33,1).Value2 = $tmpA
$wsLP1912.Cells.Item(33,3).Value2 = $tmpC
$wsLP1912.Cells.Item(33,4).Value2 = $tmpD

# swap rows 40 <-> 41
$tmpC = $wsLP1912.Cells.Item(40,3).Value2
$wsLP1912.Cells.Item(40,3).Value2 = $wsLP1912.Cells.Item(41,3).Value2
$wsLP1912.Cells.Item(41,3).Value2 = $tmpC

# swap rows 110 <-> 111
$tmpC = $wsLP1912.Cells.Item(110,3).Value2
$wsLP1912.Cells.Item(110,3).Value2 = $wsLP1912.Cells.Item(111,3).Value2
$wsLP1912.Cells.Item(111,3).Value2 = $tmpC

# swap rows 190 <-> 192
$tmpA = $wsLP1912.Cells.Item(190,1).Value2
$tmpC = $wsLP1912.Cells.Item(190,3).Value2
$tmpD = $wsLP1912.Cells.Item(190,4).Value2
$wsLP1912.Cells.Item(190,1).Value2 = $wsLP1912.Cells.Item(192,1).Value2
$wsLP1912.Cells.Item(190,3).Value2 = $wsLP1912.Cells.Item(192,3).Value2
$wsLP1912.Cells.Item(190,4).Value2 = $wsLP1912.Cells.Item(192,4).Value2
$wsLP1912.Cells.Item(192,1).Value2 = $tmpA
$wsLP1912.Cells.Item(192,3).Value2 = $tmpC
$wsLP1912.Cells.Item(192,4).Value2 = $tmpD

# swap rows 226 <-> 227
$tmpA = $wsLP1912.Cells.Item(226,1).Value2
$tmpC = $wsLP1912.Cells.Item(226,3).Value2
$tmpD = $wsLP1912.Cells.Item(226,4).Value2
$wsLP1912.Cells.Item(226,1).Value2 = $wsLP1912.Cells.Item(227,1).Value2
$wsLP1912.Cells.Item(226,3).Value2 = $wsLP1912.Cells.Item(227,3).Value2
$wsLP1912.Cells.Item(226,4).Value2 = $wsLP1912.Cells.Item(227,4).Value2
$wsLP1912.Cells.Item(227,1).Value2 = $tmpA
$wsLP1912.Cells.Item(227,3).Value2 = $tmpC
$wsLP1912.Cells.Item(227,4).Value2 = $tmpD

# swap rows 233 <-> 234
$tmpC = $wsLP1912.Cells.Item(233,3).Value2
$wsLP1912.Cells.Item(233,3).Value2 = $wsLP1912.Cells.Item(234,3).Value2
$wsLP1912.Cells.Item(234,3).Value2 = $tmpC

# --- Sheet LP1912: rewrite rows 259-289 (later scrape reshuffled/added entries) ---
$wsLP1912.Cells.Item(259,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(259,2).Value2 = "14:25"
$wsLP1912.Cells.Item(259,3).Value2 = "16_SANTA ANA"
$wsLP1912.Cells.Item(259,4).Value2 = 1
$wsLP1912.Cells.Item(259,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(260,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(260,2).Value2 = "14:28"
$wsLP1912.Cells.Item(260,3).Value2 = "15_ABASTO"
$wsLP1912.Cells.Item(260,4).Value2 = 29
$wsLP1912.Cells.Item(260,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(261,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(261,2).Value2 = "14:35"
$wsLP1912.Cells.Item(261,3).Value2 = "23_HERNANDEZ"
$wsLP1912.Cells.Item(261,4).Value2 = 11
$wsLP1912.Cells.Item(261,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(262,1).Value2 = "13:33:42"
$wsLP1912.Cells.Item(262,2).Value2 = "14:44"
$wsLP1912.Cells.Item(262,3).Value2 = "14_ABASTO"
$wsLP1912.Cells.Item(262,4).Value2 = 71
$wsLP1912.Cells.Item(262,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(263,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(263,2).Value2 = "14:44"
$wsLP1912.Cells.Item(263,3).Value2 = "15_ABASTO"
$wsLP1912.Cells.Item(263,4).Value2 = 20
$wsLP1912.Cells.Item(263,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(264,1).Value2 = "12:47:27"
$wsLP1912.Cells.Item(264,2).Value2 = "14:45"
$wsLP1912.Cells.Item(264,3).Value2 = "14_ABASTO"
$wsLP1912.Cells.Item(264,4).Value2 = 118
$wsLP1912.Cells.Item(264,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(265,1).Value2 = "12:59:47"
$wsLP1912.Cells.Item(265,2).Value2 = "14:56"
$wsLP1912.Cells.Item(265,3).Value2 = "16_P MOR-SANTA ANA"
$wsLP1912.Cells.Item(265,4).Value2 = 117
$wsLP1912.Cells.Item(265,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(266,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(266,2).Value2 = "14:57"
$wsLP1912.Cells.Item(266,3).Value2 = "16_P MOR-SANTA ANA"
$wsLP1912.Cells.Item(266,4).Value2 = 58
$wsLP1912.Cells.Item(266,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(267,1).Value2 = "12:59:47"
$wsLP1912.Cells.Item(267,2).Value2 = "14:58"
$wsLP1912.Cells.Item(267,3).Value2 = "215B_EL PATO"
$wsLP1912.Cells.Item(267,4).Value2 = 119
$wsLP1912.Cells.Item(267,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(268,1).Value2 = "13:33:42"
$wsLP1912.Cells.Item(268,2).Value2 = "15:00"
$wsLP1912.Cells.Item(268,3).Value2 = "81_EL PELIGRO"
$wsLP1912.Cells.Item(268,4).Value2 = 87
$wsLP1912.Cells.Item(268,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(269,1).Value2 = "13:33:42"
$wsLP1912.Cells.Item(269,2).Value2 = "15:05"
$wsLP1912.Cells.Item(269,3).Value2 = "10_OLMOS"
$wsLP1912.Cells.Item(269,4).Value2 = 92
$wsLP1912.Cells.Item(269,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(270,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(270,2).Value2 = "15:10"
$wsLP1912.Cells.Item(270,3).Value2 = "17_ROMERO"
$wsLP1912.Cells.Item(270,4).Value2 = 71
$wsLP1912.Cells.Item(270,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(271,1).Value2 = "13:33:42"
$wsLP1912.Cells.Item(271,2).Value2 = "15:13"
$wsLP1912.Cells.Item(271,3).Value2 = "11_ETCHEVERRY"
$wsLP1912.Cells.Item(271,4).Value2 = 100
$wsLP1912.Cells.Item(271,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(272,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(272,2).Value2 = "15:14"
$wsLP1912.Cells.Item(272,3).Value2 = "11_ETCHEVERRY"
$wsLP1912.Cells.Item(272,4).Value2 = 75
$wsLP1912.Cells.Item(272,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(273,1).Value2 = "13:33:42"
$wsLP1912.Cells.Item(273,2).Value2 = "15:17"
$wsLP1912.Cells.Item(273,3).Value2 = "26_HERNANDEZ"
$wsLP1912.Cells.Item(273,4).Value2 = 104
$wsLP1912.Cells.Item(273,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(274,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(274,2).Value2 = "15:18"
$wsLP1912.Cells.Item(274,3).Value2 = "26_HERNANDEZ"
$wsLP1912.Cells.Item(274,4).Value2 = 79
$wsLP1912.Cells.Item(274,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(275,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(275,2).Value2 = "15:21"
$wsLP1912.Cells.Item(275,3).Value2 = "26_HERNANDEZ"
$wsLP1912.Cells.Item(275,4).Value2 = 57
$wsLP1912.Cells.Item(275,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(276,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(276,2).Value2 = "15:32"
$wsLP1912.Cells.Item(276,3).Value2 = "84_COLONIA URQUIZA-ESC 49"
$wsLP1912.Cells.Item(276,4).Value2 = 68
$wsLP1912.Cells.Item(276,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(277,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(277,2).Value2 = "15:35"
$wsLP1912.Cells.Item(277,3).Value2 = "23_HERNANDEZ"
$wsLP1912.Cells.Item(277,4).Value2 = 96
$wsLP1912.Cells.Item(277,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(278,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(278,2).Value2 = "15:37"
$wsLP1912.Cells.Item(278,3).Value2 = "10_OLMOS"
$wsLP1912.Cells.Item(278,4).Value2 = 98
$wsLP1912.Cells.Item(278,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(279,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(279,2).Value2 = "15:38"
$wsLP1912.Cells.Item(279,3).Value2 = "23_HERNANDEZ"
$wsLP1912.Cells.Item(279,4).Value2 = 74
$wsLP1912.Cells.Item(279,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(280,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(280,2).Value2 = "15:39"
$wsLP1912.Cells.Item(280,3).Value2 = "215A_EL PATO"
$wsLP1912.Cells.Item(280,4).Value2 = 100
$wsLP1912.Cells.Item(280,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(281,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(281,2).Value2 = "15:46"
$wsLP1912.Cells.Item(281,3).Value2 = "14_ABASTO"
$wsLP1912.Cells.Item(281,4).Value2 = 82
$wsLP1912.Cells.Item(281,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(282,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(282,2).Value2 = "15:47"
$wsLP1912.Cells.Item(282,3).Value2 = "16_P MOR-167 Y 521"
$wsLP1912.Cells.Item(282,4).Value2 = 108
$wsLP1912.Cells.Item(282,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(283,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(283,2).Value2 = "15:48"
$wsLP1912.Cells.Item(283,3).Value2 = "14_ABASTO"
$wsLP1912.Cells.Item(283,4).Value2 = 109
$wsLP1912.Cells.Item(283,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(284,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(284,2).Value2 = "15:54"
$wsLP1912.Cells.Item(284,3).Value2 = "11_ETCHEVERRY"
$wsLP1912.Cells.Item(284,4).Value2 = 115
$wsLP1912.Cells.Item(284,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(285,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(285,2).Value2 = "15:56"
$wsLP1912.Cells.Item(285,3).Value2 = "17_ROMERO"
$wsLP1912.Cells.Item(285,4).Value2 = 92
$wsLP1912.Cells.Item(285,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(286,1).Value2 = "13:59:06"
$wsLP1912.Cells.Item(286,2).Value2 = "15:57"
$wsLP1912.Cells.Item(286,3).Value2 = "27_EL RETIRO"
$wsLP1912.Cells.Item(286,4).Value2 = 118
$wsLP1912.Cells.Item(286,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(287,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(287,2).Value2 = "16:15"
$wsLP1912.Cells.Item(287,3).Value2 = "225_C ROCA-H SUR"
$wsLP1912.Cells.Item(287,4).Value2 = 111
$wsLP1912.Cells.Item(287,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(288,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(288,2).Value2 = "16:20"
$wsLP1912.Cells.Item(288,3).Value2 = "215C_EL PATO"
$wsLP1912.Cells.Item(288,4).Value2 = 116
$wsLP1912.Cells.Item(288,5).Value2 = "LP1912"
$wsLP1912.Cells.Item(289,1).Value2 = "14:24:16"
$wsLP1912.Cells.Item(289,2).Value2 = "16:21"
$wsLP1912.Cells.Item(289,3).Value2 = "26_HERNANDEZ"
$wsLP1912.Cells.Item(289,4).Value2 = 117
$wsLP1912.Cells.Item(289,5).Value2 = "LP1912"

# --- Sheet LP1912-215: header metadata + new row 37 ---
$wsLP215.Range("A2").Value2 = "Última actualización: 14:24:16"
$wsLP215.Range("A3").Value2 = "Total filas: 32"
$wsLP215.Cells.Item(37,1).Value2 = "14:24:16"
$wsLP215.Cells.Item(37,2).Value2 = "16:20"
$wsLP215.Cells.Item(37,3).Value2 = "215C_EL PATO"
$wsLP215.Cells.Item(37,4).Value2 = 116
$wsLP215.Cells.Item(37,5).Value2 = "LP1912"

# --- Sheet 6203-6173: header metadata + new row 47 ---
$wsL6203.Range("A2").Value2 = "Última actualización: 14:24:16"
$wsL6203.Range("A3").Value2 = "Total filas: 42"
$wsL6203.Cells.Item(47,1).Value2 = "14:24:16"
$wsL6203.Cells.Item(47,2).Value2 = "16:14"
$wsL6203.Cells.Item(47,3).Value2 = "215C_LA PLATA"
$wsL6203.Cells.Item(47,4).Value2 = 110
$wsL6203.Cells.Item(47,5).Value2 = "L6203"

Write-Output "done"
